$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated cell values from the diff.
# Column D sometimes holds numeric-looking text (e.g. "6.160", "24.704.10") that must
# stay literal text, so each such cell is pre-formatted as Text before the value is set.

$ws.Range("D2").Value = "24.704.10"
$ws.Range("E2").Value = "  -2.26%  "
$ws.Range("D3").Value = "1.673.81"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +1.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.78"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3673"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3351"
$ws.Range("E8").Value = "  -4.52%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.54"
$ws.Range("E9").Value = "  -5.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.167"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07297"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.160"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.44"
$ws.Range("E14").Value = "  -3.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.798"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "1.676.77"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001097"
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9972"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.72"
$ws.Range("E20").Value = "  -3.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.84"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.161"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.55"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("D24").Value = "24.725.96"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.703"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.76"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.83"
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "129.65"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "1.865.37"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.220"
$ws.Range("E31").Value = "  +10.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.489"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.158"
$ws.Range("E33").Value = "  +3.36%  "
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.750"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.37"
$ws.Range("E35").Value = "  +6.05%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08593"
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.421"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06462"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.794"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2153"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.242"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6239"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9958"
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.42"
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.785"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5948"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.045"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.07"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07134"
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.76"
$ws.Range("E51").Value = "  -0.14%  "
